$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to remain a text value even if it looks numeric,
# mirroring the source workbook where these columns are stored as text.
function Set-TextCell($rng, [string]$val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '27.014.50'
$ws.Range("E2").Value = '  +1.93%  '

$ws.Range("D3").Value = '1.674.13'
$ws.Range("E3").Value = '  +3.04%  '

$ws.Range("E4").Value = '  +0.12%  '

Set-TextCell $ws.Range("D5") '216.59'
$ws.Range("E5").Value = '  +1.47%  '

Set-TextCell $ws.Range("D6") '0.535'
$ws.Range("E6").Value = '  +6.67%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("E8").Value = '  +2.80%  '

Set-TextCell $ws.Range("D9") '0.0621'
$ws.Range("E9").Value = '  +1.93%  '

Set-TextCell $ws.Range("D10") '20.22'
$ws.Range("E10").Value = '  +5.20%  '

Set-TextCell $ws.Range("D11") '0.0892'
$ws.Range("E11").Value = '  +4.36%  '

$ws.Range("D12").Value = '1.910.20'
$ws.Range("E12").Value = '  +3.11%  '

$ws.Range("D13").Value = '1.674.79'
$ws.Range("E13").Value = '  +3.14%  '

$ws.Range("E14").Value = '  +1.16%  '

Set-TextCell $ws.Range("D15") '65.83'
$ws.Range("E15").Value = '  +2.93%  '

Set-TextCell $ws.Range("D16") '0.521'
$ws.Range("E16").Value = '  +1.93%  '

$ws.Range("D17").Value = '27.032.29'

Set-TextCell $ws.Range("D18") '233.91'
$ws.Range("E18").Value = '  -0.13%  '

$ws.Range("D19").Value = '0.0₃0738'
$ws.Range("E19").Value = '  +1.77%  '

Set-TextCell $ws.Range("D20") '7.77'
$ws.Range("E20").Value = '  +0.60%  '

$ws.Range("E21").Value = '  -0.04%  '

Set-TextCell $ws.Range("D22") '4.47'
$ws.Range("E22").Value = '  +3.45%  '

$ws.Range("E23").Value = '  +0.46%  '

Set-TextCell $ws.Range("D24") '9.27'
$ws.Range("E24").Value = '  +1.28%  '

Set-TextCell $ws.Range("D25") '145.89'
$ws.Range("E25").Value = '  -0.08%  '

Set-TextCell $ws.Range("D26") '0.117'
$ws.Range("E26").Value = '  +3.11%  '

$ws.Range("E27").Value = '  +1.15%  '

Set-TextCell $ws.Range("D28") '15.89'
$ws.Range("E28").Value = '  +1.71%  '

$ws.Range("E29").Value = '  +0.09%  '

Set-TextCell $ws.Range("D30") '0.0498'
$ws.Range("E30").Value = '  +0.89%  '

$ws.Range("E31").Value = '  +1.63%  '

$ws.Range("E32").Value = '  +1.80%  '

$ws.Range("D33").Value = '1.453.39'
$ws.Range("E33").Value = '  -4.36%  '

Set-TextCell $ws.Range("D34") '3.17'
$ws.Range("E34").Value = '  +5.88%  '

$ws.Range("E35").Value = '  +5.97%  '

$ws.Range("E36").Value = '  -0.52%  '

Set-TextCell $ws.Range("D37") '0.899'

$ws.Range("E38").Value = '  -0.69%  '

Set-TextCell $ws.Range("D39") '0.0169'
$ws.Range("E39").Value = '  +1.65%  '

Set-TextCell $ws.Range("D40") '6.06'
$ws.Range("E40").Value = '  +3.63%  '

$ws.Range("E41").Value = '  +0.04%  '

$ws.Range("E42").Value = '  +4.30%  '

$ws.Range("E43").Value = '  +7.56%  '

Set-TextCell $ws.Range("D44") '65.91'
$ws.Range("E44").Value = '  +5.19%  '

$ws.Range("D45").Value = '1.816.56'
$ws.Range("E45").Value = '  +2.97%  '

Set-TextCell $ws.Range("D46") '0.785'
$ws.Range("E46").Value = '  +3.22%  '

Set-TextCell $ws.Range("D47") '90.62'
$ws.Range("E47").Value = '  +0.82%  '

$ws.Range("E48").Value = '  +1.46%  '

Set-TextCell $ws.Range("D49") '0.101'
$ws.Range("E49").Value = '  +4.61%  '

Set-TextCell $ws.Range("D50") '0.0508'
$ws.Range("E50").Value = '  +1.21%  '

Set-TextCell $ws.Range("D51") '7.65'
$ws.Range("E51").Value = '  +1.30%  '
